$wb = $excel.ActiveWorkbook

$wsAbout      = $wb.Worksheets.Item("About")
$wsInteger    = $wb.Worksheets.Item("Integer")
$wsBoolean    = $wb.Worksheets.Item("Boolean")
$wsSubscript  = $wb.Worksheets.Item("Subscript")

# ------------------------------------------------------------------
# About sheet - note text updates (shared string reshuffle, no visible
# text change other than the note that used to be the last one).
# ------------------------------------------------------------------
$wsAbout.Range("A7").Value = "InputData pathnames of CSV files with values constrained to specific data types"

# ------------------------------------------------------------------
# Integer sheet - content unchanged, values restated for safety.
# ------------------------------------------------------------------
$wsInteger.Range("A1").Value = "InputData pathname"
$wsInteger.Range("A2").Value = "bldgs/CL/CL.csv"
$wsInteger.Range("A3").Value = "ccs/BCS/BCS-DoSfCS.csv"
$wsInteger.Range("A4").Value = "elec/DRC/DRC-ADRHpDRE.csv"
$wsInteger.Range("A5").Value = "elec/GBSC/GBDSD.csv"
$wsInteger.Range("A6").Value = "elec/MLfPPR/MLfPPR.csv"
$wsInteger.Range("A7").Value = "elec/RPfFESCC/RPfFESCC.csv"
$wsInteger.Range("A8").Value = "fuels/BS/BS-DoSpUEO.csv"
$wsInteger.Range("A9").Value = "trans/AVL/AVL.csv"

# ------------------------------------------------------------------
# Boolean sheet - the "trans/BVTQaZ/BVTQaZ.csv" and "trans/VTQaZ/VTQaZ.csv"
# single-file rows are each split into six per-vehicle-type files.
# ------------------------------------------------------------------
$wsBoolean.Range("A1").Value  = "InputData pathname"
$wsBoolean.Range("A2").Value  = "ctrl-settings/BAEPAbCiPC/BAEPAbCiPC.csv"
$wsBoolean.Range("A3").Value  = "ctrl-settings/BDCTBA/BDCTBA.csv"
$wsBoolean.Range("A4").Value  = "ctrl-settings/BDMFL/BDMFL.csv"
$wsBoolean.Range("A5").Value  = "ctrl-settings/BENCEfCT/BENCEfCT.csv"
$wsBoolean.Range("A6").Value  = "ctrl-settings/BEPEfCT/BEPEfCT.csv"
$wsBoolean.Range("A7").Value  = "ctrl-settings/BESHFoFRV/BESHFoFRV.csv"
$wsBoolean.Range("A8").Value  = "ctrl-settings/BIEfEE/BIEfEE.csv"
$wsBoolean.Range("A9").Value  = "ctrl-settings/BIEfIE/BIEfIE.csv"
$wsBoolean.Range("A10").Value = "ctrl-settings/BRCToEP/BRCToEP.csv"
$wsBoolean.Range("A11").Value = "ctrl-settings/BUTYGV/BUTYGV.csv"
$wsBoolean.Range("A12").Value = "ctrl-settings/EGGRA/EGGRA-use-adjustment.csv"
$wsBoolean.Range("A13").Value = "elec/BDTPTUMCF/BDTPTUMCF.csv"
$wsBoolean.Range("A14").Value = "elec/RQSD/RQSD-BRQSD.csv"
$wsBoolean.Range("A15").Value = "elec/RQSD/RQSD-RQSD.csv"
$wsBoolean.Range("A16").Value = "indst/IFStFS/IFStFS.csv"

$wsBoolean.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBoolean.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBoolean.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBoolean.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBoolean.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBoolean.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

$wsBoolean.Range("A23").Value = "trans/BVTStL/BVTStL.csv"
$wsBoolean.Range("A24").Value = "trans/PVTStL/PVTStL.csv"
$wsBoolean.Range("A25").Value = "trans/SRPbVT/SRPbVT.csv"

$wsBoolean.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBoolean.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBoolean.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBoolean.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBoolean.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBoolean.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

$wsBoolean.Range("A32").Value = "trans/VTStFES/VTStFES.csv"

# Trailing formatted-but-empty rows (33-38) that exist below the data in
# the final workbook. Touch their font so the rows persist with the same
# cell style used throughout column A, and the sheet dimension grows to
# A1:A38 to match.
for ($r = 33; $r -le 38; $r++) {
    $wsBoolean.Rows.Item($r).Font.Name = "Calibri"
    $wsBoolean.Rows.Item($r).Font.Size = 11
}

# ------------------------------------------------------------------
# Subscript sheet - content unchanged, values restated for safety.
# ------------------------------------------------------------------
$wsSubscript.Range("A1").Value = "InputData pathname"
$wsSubscript.Range("A2").Value = "elec/ESUfRaLCD/ESUfRaLCD-dispatch.csv"
$wsSubscript.Range("A3").Value = "elec/ESUfRaLCD/ESUfRaLCD-reliability.csv"
$wsSubscript.Range("A4").Value = "plcy-schd/FoPITY/FoPITY-policy-elements.csv"
$wsSubscript.Range("A5").Value = "plcy-schd/FY/FY.csv"

# ------------------------------------------------------------------
# Window / selection state.
# Integer ends up with the cursor parked at A13, Boolean at A32 (after
# scrolling so row 10 is near the top), and About becomes the active
# (selected) tab.
# ------------------------------------------------------------------
$wsInteger.Activate()
$wsInteger.Range("A13").Select()

$wsBoolean.Activate()
$wsBoolean.Range("A32").Select()

$wsAbout.Activate()
$wsAbout.Range("A1").Select()
